{"js": "// Remove the inline picture (Screenshot 2024-09-25 144741.png) that was\n// uploaded into the document, leaving its (now empty) paragraph intact.\nconst body = context.document.body;\nconst pics = body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nfor (let i = pics.items.length - 1; i >= 0; i--) {\n  pics.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the inline picture (Screenshot 2024-09-25 144741.png) that was\n# uploaded into the document, leaving its (now empty) paragraph intact.\n$d = $word.ActiveDocument\n\nfor ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {\n    $d.InlineShapes($i).Delete()\n}\n"}
